# Automatische test-sync: 2025-07-31 21:25:50
#
# Adds a second "Testmail #2" row to the Logs sheet, adds the matching
# rollup row to the Dashboard sheet, grows the chart's category/value
# series references to include the new Dashboard row, and widens the
# Logs sheet's conditional-formatting ranges to cover the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append row 4 with the new test-mail data
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Wil je dit oppakken?"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Testmail #2: Wil je dit oppakken?"
$logs.Range("D4").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E4").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$logs.Range("F4").Value = "2025-07-31 21:25:23"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# ---------------------------------------------------------------
# 2. Logs sheet: widen the conditional-formatting ranges from
#    row 2:3 to row 2:4 (D, G, H, I, J columns)
# ---------------------------------------------------------------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`3")
    $newRange = $logs.Range("$col`2:$col`4")
    $rules = $oldRange.FormatConditions
    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------
# 3. Dashboard sheet: append row 3 rollup for the new category
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------
# 4. Chart: extend the category/value series to A2:A3 / B2:B3
# ---------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$3,'Dashboard'!`$B`$2:`$B`$3,1)"
